$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.294.49'
$ws.Range("E2").Value = '  +0.78%  '
$ws.Range("D3").Value = '3.767.59'
$ws.Range("E3").Value = '  +0.54%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '620.30'
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.36'
$ws.Range("E6").Value = '  +2.12%  '
$ws.Range("D7").Value = '3.765.56'
$ws.Range("E7").Value = '  +0.55%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.536'
$ws.Range("E9").Value = '  -1.52%  '
$ws.Range("E10").Value = '  -0.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.51'
$ws.Range("E11").Value = '  +1.88%  '
$ws.Range("E12").Value = '  -3.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.21'
$ws.Range("E13").Value = '  -2.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000258'
$ws.Range("E14").Value = '  +0.37%  '
$ws.Range("D15").Value = '4.407.20'
$ws.Range("E15").Value = '  +0.77%  '
$ws.Range("D16").Value = '3.767.41'
$ws.Range("E16").Value = '  +0.46%  '
$ws.Range("D17").Value = '70.317.28'
$ws.Range("E17").Value = '  +0.59%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.61'
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("E19").Value = '  -2.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.68'
$ws.Range("E20").Value = '  -0.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '507.58'
$ws.Range("E21").Value = '  -2.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.28'
$ws.Range("E22").Value = '  -1.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.729'
$ws.Range("E23").Value = '  -1.53%  '
$ws.Range("E24").Value = '  +5.59%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.95'
$ws.Range("E25").Value = '  -2.37%  '
$ws.Range("B26").Value = 'RenderToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.40'
$ws.Range("E26").Value = '  +4.29%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '13.09'
$ws.Range("E27").Value = '  -3.81%  '
$ws.Range("E28").Value = '  +7.00%  '
$ws.Range("E29").Value = '  +0.34%  '
$ws.Range("E30").Value = '  -0.12%  '
$ws.Range("E31").Value = '  +2.59%  '
$ws.Range("E32").Value = '  +2.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.91'
$ws.Range("E33").Value = '  -2.50%  '
$ws.Range("E34").Value = '  -1.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.26%  '
$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.22'
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("B37").Value = 'Mantle'
$ws.Range("C37").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.06'
$ws.Range("E37").Value = '  +1.51%  '
$ws.Range("E38").Value = '  +2.47%  '
$ws.Range("E39").Value = '  +6.13%  '
$ws.Range("E40").Value = '  +13.05%  '
$ws.Range("E41").Value = '  -4.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '45.91'
$ws.Range("E42").Value = '  +2.91%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '50.03'
$ws.Range("E43").Value = '  -2.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '440.08'
$ws.Range("E44").Value = '  +2.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.68'
$ws.Range("E45").Value = '  -2.14%  '
$ws.Range("D46").Value = '2.991.90'
$ws.Range("E46").Value = '  -2.89%  '
$ws.Range("E47").Value = '  +0.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.67'
$ws.Range("E48").Value = '  -1.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '138.88'
$ws.Range("E49").Value = '  +1.72%  '
$ws.Range("E50").Value = '  -0.03%  '
$ws.Range("E51").Value = '  -0.80%  '
